$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1", "mathbr", "2000-01-01", "P"),
    @("1", "mathbr", "2000-01-02", "P"),
    @("1", "mathbr", "2000-01-03", "A"),
    @("1", "mathbr", "2000-01-04", "P"),
    @("1", "mathbr", "2000-01-05", "A"),
    @("1", "mathbr", "2000-01-07", "P")
)

$startRow = 60
$rng = $ws.Range("A${startRow}:D65")
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}
